$d = $word.ActiveDocument

# --- Paragraph 2: the "{m:userdoc 'zone1'}" field -> literal braces text ---
$p1 = $d.Paragraphs.Item(2)
$r1 = $p1.Range
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' `
  + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
  + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
  + '<pkg:xmlData>' `
  + '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' `
  + '<w:body>' `
  + '<w:p>' `
  + '<w:r><w:t>{</w:t></w:r>' `
  + '<w:r><w:t>m</w:t></w:r>' `
  + "<w:r><w:t>:userdoc 'zone1'</w:t></w:r>" `
  + '<w:r><w:t xml:space="preserve">}</w:t></w:r>' `
  + '</w:p>' `
  + '</w:body>' `
  + '</w:document>' `
  + '</pkg:xmlData></pkg:part></pkg:package>'
$r1.InsertXML($xml1)

# --- Paragraph 4: the "{m:enduserdoc}" field (with _GoBack bookmark) -> literal braces text ---
$p2 = $d.Paragraphs.Item(4)
$r2 = $p2.Range
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' `
  + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
  + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
  + '<pkg:xmlData>' `
  + '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' `
  + '<w:body>' `
  + '<w:p>' `
  + '<w:r><w:t>{m:</w:t></w:r>' `
  + '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' `
  + '<w:bookmarkEnd w:id="0"/>' `
  + '<w:r><w:t xml:space="preserve">enduserdoc}</w:t></w:r>' `
  + '</w:p>' `
  + '</w:body>' `
  + '</w:document>' `
  + '</pkg:xmlData></pkg:part></pkg:package>'
$r2.InsertXML($xml2)

Write-Output "edit applied"
